# "get dhamma.org 'Other' from other_course in central db"
#
# The "Other" / "Trust WE" / "|TRUST MEETING|" row (row 5) is no longer
# needed because that mapping is now looked up straight from the central
# db's other_course table, so the whole row is removed. That also empties
# out column C ("course_type") everywhere else in the sheet, so the column
# is dropped entirely too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Other" row (row 5) and shift the rows below it up.
$ws.Rows("5:5").Delete()

# Column C ("course_type") is now unused - remove it and shift left.
$ws.Columns("C:C").Delete()

# Leave the selection where the old column C header used to be, matching
# the author's final cursor position.
[void]$ws.Range("C1").Select()
